$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 58
$ws.Range("H58").Value = 425.1111
$ws.Range("I58").Value = 313.7143
$ws.Range("J58").Value = 815
$ws.Range("K58").Value = 941.1428999999999
$ws.Range("L58").Value = 2445
$ws.Range("M58").Value = -791.1428999999999
$ws.Range("N58").Value = -2745

# Row 61
$ws.Range("H61").Value = 55555
$ws.Range("I61").Value = 55555
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 166665
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -166493

# Row 113
$ws.Range("H113").Value = 8549.666999999999
$ws.Range("I113").Value = 7766.3335
$ws.Range("J113").Value = 9333
$ws.Range("K113").Value = 7766.3335
$ws.Range("L113").Value = 9333
$ws.Range("M113").Value = -4512.3335
$ws.Range("N113").Value = -15841

# Row 129
$ws.Range("H129").Value = 1539.1666
$ws.Range("I129").Value = 560
$ws.Range("J129").Value = 3497.5
$ws.Range("K129").Value = 1680
$ws.Range("L129").Value = 10492.5
$ws.Range("M129").Value = 3320
$ws.Range("N129").Value = -20492.5

# Row 132
$ws.Range("H132").Value = 1807.9
$ws.Range("I132").Value = 1785
$ws.Range("J132").Value = 1899.5
$ws.Range("K132").Value = 5355
$ws.Range("L132").Value = 5698.5
$ws.Range("M132").Value = -2825
$ws.Range("N132").Value = -10758.5

# Row 135
$ws.Range("H135").Value = 899.25
$ws.Range("I135").Value = 884.8570999999999
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 7963.7139
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -5428.7139
$ws.Range("N135").Value = -14070

# Row 141
$ws.Range("H141").Value = 2272.875
$ws.Range("I141").Value = 2312.1428
$ws.Range("J141").Value = 1998
$ws.Range("K141").Value = 6936.428400000001
$ws.Range("L141").Value = 5994
$ws.Range("M141").Value = -1756.428400000001
$ws.Range("N141").Value = -16354

$ws = $wb.Worksheets.Item("ARM")
# Row 19
$ws.Range("H19").Value = 5266.3335
$ws.Range("I19").Value = 2899.5
$ws.Range("J19").Value = 10000
$ws.Range("K19").Value = 2899.5
$ws.Range("L19").Value = 10000
$ws.Range("M19").Value = -2670.5
$ws.Range("N19").Value = -10458

# Row 32
$ws.Range("H32").Value = 1391.7255
$ws.Range("I32").Value = 1040.3877
$ws.Range("J32").Value = 9999.5
$ws.Range("K32").Value = 1040.3877
$ws.Range("L32").Value = 9999.5
$ws.Range("M32").Value = -753.3877
$ws.Range("N32").Value = -10573.5

# Row 61
$ws.Range("H61").Value = 2730.5
$ws.Range("I61").Value = 2307.8333
$ws.Range("J61").Value = 3998.5
$ws.Range("K61").Value = 2307.8333
$ws.Range("L61").Value = 3998.5
$ws.Range("M61").Value = -2095.8333
$ws.Range("N61").Value = -4422.5

# Row 136
$ws.Range("H136").Value = 2730.5
$ws.Range("I136").Value = 2307.8333
$ws.Range("J136").Value = 3998.5
$ws.Range("K136").Value = 6923.499899999999
$ws.Range("L136").Value = 11995.5
$ws.Range("M136").Value = -4373.499899999999
$ws.Range("N136").Value = -17095.5

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1660
$ws.Range("I94").Value = 1690.6666
$ws.Range("J94").Value = 1200
$ws.Range("K94").Value = 1690.6666
$ws.Range("L94").Value = 1200
$ws.Range("M94").Value = -1239.6666
$ws.Range("N94").Value = -2102

$ws = $wb.Worksheets.Item("CRP")
# Row 14
$ws.Range("H14").Value = 7405
$ws.Range("I14").Value = 799
$ws.Range("J14").Value = 14011
$ws.Range("K14").Value = 799
$ws.Range("L14").Value = 14011
$ws.Range("M14").Value = -629
$ws.Range("N14").Value = -14351

# Row 16
$ws.Range("H16").Value = 536.8333
$ws.Range("I16").Value = 536.8333
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 536.8333
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -249.8333

# Row 31
$ws.Range("H31").Value = 4525.6665
$ws.Range("I31").Value = 2797.5
$ws.Range("J31").Value = 5389.75
$ws.Range("K31").Value = 2797.5
$ws.Range("L31").Value = 5389.75
$ws.Range("M31").Value = -2502.5
$ws.Range("N31").Value = -5979.75

# Row 34
$ws.Range("H34").Value = 4525.6665
$ws.Range("I34").Value = 2797.5
$ws.Range("J34").Value = 5389.75
$ws.Range("K34").Value = 2797.5
$ws.Range("L34").Value = 5389.75
$ws.Range("M34").Value = -2595.5
$ws.Range("N34").Value = -5793.75

# Row 50
$ws.Range("H50").Value = 20084
$ws.Range("I50").Value = 20126.625
$ws.Range("J50").Value = 19998.75
$ws.Range("K50").Value = 20126.625
$ws.Range("L50").Value = 19998.75
$ws.Range("M50").Value = -19501.625
$ws.Range("N50").Value = -21248.75

# Row 96
$ws.Range("H96").Value = 21655
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 21655
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 21655
$ws.Range("N96").Value = -27147

# Row 113
$ws.Range("H113").Value = 536.8333
$ws.Range("I113").Value = 536.8333
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 536.8333
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1633.1667

# Row 132
$ws.Range("H132").Value = 2405
$ws.Range("I132").Value = 2670.3333
$ws.Range("J132").Value = 2007
$ws.Range("K132").Value = 8010.999899999999
$ws.Range("L132").Value = 6021
$ws.Range("M132").Value = -5480.999899999999
$ws.Range("N132").Value = -11081

# Row 134
$ws.Range("H134").Value = 1000.7
$ws.Range("I134").Value = 963.375
$ws.Range("J134").Value = 1150
$ws.Range("K134").Value = 2890.125
$ws.Range("L134").Value = 3450
$ws.Range("M134").Value = -355.125
$ws.Range("N134").Value = -8520

$ws = $wb.Worksheets.Item("CUL")
# Row 6
$ws.Range("H6").Value = 3440.1
$ws.Range("I6").Value = 489
$ws.Range("J6").Value = 30000
$ws.Range("K6").Value = 1467
$ws.Range("L6").Value = 90000
$ws.Range("M6").Value = -1354
$ws.Range("N6").Value = -90226

# Row 68
$ws.Range("H68").Value = 1499.3334
$ws.Range("I68").Value = 1499
$ws.Range("J68").Value = 1500
$ws.Range("K68").Value = 4497
$ws.Range("L68").Value = 4500
$ws.Range("M68").Value = -3686
$ws.Range("N68").Value = -6122

# Row 71
$ws.Range("H71").Value = 1499.3334
$ws.Range("I71").Value = 1499
$ws.Range("J71").Value = 1500
$ws.Range("K71").Value = 13491
$ws.Range("L71").Value = 13500
$ws.Range("M71").Value = -9435
$ws.Range("N71").Value = -21612

# Row 113
$ws.Range("H113").Value = 1384.25
$ws.Range("I113").Value = 337.5
$ws.Range("J113").Value = 1733.1666
$ws.Range("K113").Value = 1012.5
$ws.Range("L113").Value = 5199.4998
$ws.Range("M113").Value = 1157.5
$ws.Range("N113").Value = -9539.4998

$ws = $wb.Worksheets.Item("GSM")
# Row 24
$ws.Range("H24").Value = 35000000
$ws.Range("I24").Value = 35000000
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 35000000
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -34999827
$ws.Range("N24").ClearContents()

# Row 102
$ws.Range("H102").Value = 5735.4
$ws.Range("I102").Value = 5735.4
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 5735.4
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -4113.4

# Row 135
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 31276.908
$ws.Range("I7").Value = 33256.145
$ws.Range("J7").Value = 27813.25
$ws.Range("K7").Value = 33256.145
$ws.Range("L7").Value = 27813.25
$ws.Range("M7").Value = -33144.145
$ws.Range("N7").Value = -28037.25

# Row 22
$ws.Range("H22").Value = 957.8182
$ws.Range("I22").Value = 889.8333
$ws.Range("J22").Value = 1039.4
$ws.Range("K22").Value = 889.8333
$ws.Range("L22").Value = 1039.4
$ws.Range("M22").Value = -594.8333
$ws.Range("N22").Value = -1629.4

# Row 27
$ws.Range("H27").Value = 957.8182
$ws.Range("I27").Value = 889.8333
$ws.Range("J27").Value = 1039.4
$ws.Range("K27").Value = 889.8333
$ws.Range("L27").Value = 1039.4
$ws.Range("M27").Value = -782.8333
$ws.Range("N27").Value = -1253.4

# Row 40
$ws.Range("H40").Value = 5543.3335
$ws.Range("I40").Value = 5298.75
$ws.Range("J40").Value = 7500
$ws.Range("K40").Value = 5298.75
$ws.Range("L40").Value = 7500
$ws.Range("M40").Value = -5162.75
$ws.Range("N40").Value = -7772

# Row 75
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()

# Row 78
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()

# Row 81
$ws.Range("H81").Value = 50000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 50000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996

# Row 82
$ws.Range("H82").Value = 1779.4
$ws.Range("I82").Value = 1474.25
$ws.Range("J82").Value = 3000
$ws.Range("K82").Value = 1474.25
$ws.Range("L82").Value = 3000
$ws.Range("M82").Value = -1113.25
$ws.Range("N82").Value = -3722

# Row 84
$ws.Range("H84").Value = 50000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 50000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984

# Row 85
$ws.Range("H85").Value = 1779.4
$ws.Range("I85").Value = 1474.25
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 1474.25
$ws.Range("L85").Value = 3000
$ws.Range("M85").Value = -226.25
$ws.Range("N85").Value = -5496

# Row 126
$ws.Range("H126").Value = 31276.908
$ws.Range("I126").Value = 33256.145
$ws.Range("J126").Value = 27813.25
$ws.Range("K126").Value = 99768.435
$ws.Range("L126").Value = 83439.75
$ws.Range("M126").Value = -97298.435
$ws.Range("N126").Value = -88379.75

# Row 139
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 100278.125
$ws.Range("I5").Value = 125000
$ws.Range("J5").Value = 26112.5
$ws.Range("K5").Value = 125000
$ws.Range("L5").Value = 26112.5
$ws.Range("M5").Value = -124888
$ws.Range("N5").Value = -26336.5

# Row 126
$ws.Range("H126").Value = 99999
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 99999
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 299997
$ws.Range("N126").Value = -304937
$ws.Range("M126").ClearContents()

# Row 132
$ws.Range("H132").Value = 2529.3635
$ws.Range("I132").Value = 1980.4445
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 5941.333500000001
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -3411.333500000001
$ws.Range("N132").Value = -20058.5

# Row 136
$ws.Range("H136").Value = 2469.7058
$ws.Range("I136").Value = 3407.818
$ws.Range("J136").Value = 749.8333
$ws.Range("K136").Value = 10223.454
$ws.Range("L136").Value = 2249.4999
$ws.Range("M136").Value = -7673.454000000002
$ws.Range("N136").Value = -7349.4999
